$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same cell style (horizontal+vertical center) used throughout the
# existing table to the new row range.
$ws.Range("A178:J207").HorizontalAlignment = -4108
$ws.Range("A178:J207").VerticalAlignment = -4108

$ws.Cells.Item(178,1).Value = "밀리아, 프레이야, 쥬리"
$ws.Cells.Item(178,2).Value = "파이크"
$ws.Cells.Item(178,4).Value = "밀리아, 프레이야, 바네사"
$ws.Cells.Item(178,5).Value = "루"
$ws.Cells.Item(178,6).Value = "밀2프2프1"
$ws.Cells.Item(178,3).Value = "쥬2프2프1"
$ws.Cells.Item(178,7).Value = "선"
$ws.Cells.Item(178,8).Value = 260103
$ws.Cells.Item(178,9).Value = "느그클럽"
$ws.Cells.Item(178,10).Value = "공격"
$ws.Cells.Item(179,1).Value = "플라튼, 실베스타, 아멜리아"
$ws.Cells.Item(179,2).Value = "루"
$ws.Cells.Item(179,4).Value = "카일, 엘리시아, 파이"
$ws.Cells.Item(179,5).Value = "이린"
$ws.Cells.Item(179,6).Value = "엘1카1카2"
$ws.Cells.Item(179,3).Value = "아1실2플2"
$ws.Cells.Item(179,7).Value = "후"
$ws.Cells.Item(179,8).Value = 260103
$ws.Cells.Item(179,9).Value = "느그클럽"
$ws.Cells.Item(179,10).Value = "공격"
$ws.Cells.Item(180,1).Value = "오공, 겔리두스, 엘리시아"
$ws.Cells.Item(180,2).Value = "루"
$ws.Cells.Item(180,4).Value = "오공, 스파이크, 겔리두스"
$ws.Cells.Item(180,5).Value = "유"
$ws.Cells.Item(180,6).Value = "오2겔2겔1"
$ws.Cells.Item(180,3).Value = "오2겔2엘2"
$ws.Cells.Item(180,7).Value = "선"
$ws.Cells.Item(180,8).Value = 260103
$ws.Cells.Item(180,9).Value = "느그클럽"
$ws.Cells.Item(180,10).Value = "공격"
$ws.Cells.Item(181,1).Value = "플라튼, 콜트, 실베스타"
$ws.Cells.Item(181,2).Value = "카람"
$ws.Cells.Item(181,4).Value = "연희, 바네사, 아멜리아"
$ws.Cells.Item(181,5).Value = "연지"
$ws.Cells.Item(181,6).Value = "아2연2연1"
$ws.Cells.Item(181,3).Value = "플2콜1콜2"
$ws.Cells.Item(181,7).Value = "선"
$ws.Cells.Item(181,8).Value = 260103
$ws.Cells.Item(181,9).Value = "느그클럽"
$ws.Cells.Item(181,10).Value = "공격"
$ws.Cells.Item(182,1).Value = "프레이야, 바네사, 밀리아"
$ws.Cells.Item(182,2).Value = "연지"
$ws.Cells.Item(182,4).Value = "오공, 스파이크, 겔리두스"
$ws.Cells.Item(182,5).Value = "크리"
$ws.Cells.Item(182,6).Value = "오2겔2오1"
$ws.Cells.Item(182,3).Value = "프2바1프1"
$ws.Cells.Item(182,7).Value = "후"
$ws.Cells.Item(182,8).Value = 260103
$ws.Cells.Item(182,9).Value = "느그클럽"
$ws.Cells.Item(182,10).Value = "공격"
$ws.Cells.Item(183,1).Value = "플라튼, 트루드, 아멜리아"
$ws.Cells.Item(183,2).Value = "카람"
$ws.Cells.Item(183,4).Value = "크리스, 녹스, 로지"
$ws.Cells.Item(183,5).Value = "맬패로"
$ws.Cells.Item(183,6).Value = "크1크2녹2"
$ws.Cells.Item(183,3).Value = "아1트2트1"
$ws.Cells.Item(183,7).Value = "선"
$ws.Cells.Item(183,8).Value = 260103
$ws.Cells.Item(183,9).Value = "느그클럽"
$ws.Cells.Item(183,10).Value = "공격"
$ws.Cells.Item(184,1).Value = "겔리두스, 스파이크, 챈슬러"
$ws.Cells.Item(184,2).Value = "맬패로"
$ws.Cells.Item(184,4).Value = "플라튼, 엘리스, 라니아"
$ws.Cells.Item(184,5).Value = "루"
$ws.Cells.Item(184,6).Value = "플2라2엘2"
$ws.Cells.Item(184,3).Value = "챈2챈1"
$ws.Cells.Item(184,7).Value = "후"
$ws.Cells.Item(184,8).Value = 260103
$ws.Cells.Item(184,9).Value = "느그클럽"
$ws.Cells.Item(184,10).Value = "공격"
$ws.Cells.Item(185,1).Value = "카르마, 챈슬러, 크리스"
$ws.Cells.Item(185,2).Value = "파이크"
$ws.Cells.Item(185,4).Value = "카일, 카구라, 파이"
$ws.Cells.Item(185,5).Value = "이린"
$ws.Cells.Item(185,6).Value = "구2카1카2"
$ws.Cells.Item(185,3).Value = "크1카2"
$ws.Cells.Item(185,7).Value = "선"
$ws.Cells.Item(185,8).Value = 260103
$ws.Cells.Item(185,9).Value = "느그클럽"
$ws.Cells.Item(185,10).Value = "공격"
$ws.Cells.Item(186,1).Value = "트루드, 스파이크, 아멜리아"
$ws.Cells.Item(186,2).Value = "파이크"
$ws.Cells.Item(186,4).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(186,5).Value = "연지"
$ws.Cells.Item(186,6).Value = "바1프2키2"
$ws.Cells.Item(186,3).Value = "아2아1트2"
$ws.Cells.Item(186,7).Value = "선"
$ws.Cells.Item(186,8).Value = 260103
$ws.Cells.Item(186,9).Value = "느그클럽"
$ws.Cells.Item(186,10).Value = "공격"
$ws.Cells.Item(187,1).Value = "연희, 바네사, 로지"
$ws.Cells.Item(187,2).Value = "연지"
$ws.Cells.Item(187,4).Value = "플라튼, 실베스타, 키리엘"
$ws.Cells.Item(187,5).Value = "루"
$ws.Cells.Item(187,6).Value = "키2실2실1"
$ws.Cells.Item(187,3).Value = "바1연1연2"
$ws.Cells.Item(187,7).Value = "후"
$ws.Cells.Item(187,8).Value = 260103
$ws.Cells.Item(187,9).Value = "느그클럽"
$ws.Cells.Item(187,10).Value = "공격"
$ws.Cells.Item(188,1).Value = "프레이야, 바네사, 쥬리"
$ws.Cells.Item(188,2).Value = "노트"
$ws.Cells.Item(188,4).Value = "오공, 스파이크, 콜트"
$ws.Cells.Item(188,5).Value = "루"
$ws.Cells.Item(188,6).Value = "오2콜1콜2"
$ws.Cells.Item(188,3).Value = "프2쥬2프1"
$ws.Cells.Item(188,7).Value = "후"
$ws.Cells.Item(188,8).Value = 260103
$ws.Cells.Item(188,9).Value = "느그클럽"
$ws.Cells.Item(188,10).Value = "공격"
$ws.Cells.Item(189,1).Value = "트루드, 스파이크, 아멜리아"
$ws.Cells.Item(189,2).Value = "크리"
$ws.Cells.Item(189,4).Value = "오공, 스파이크, 아멜리아"
$ws.Cells.Item(189,5).Value = "크리"
$ws.Cells.Item(189,6).Value = "아2오2스2"
$ws.Cells.Item(189,3).Value = "아1아2트2"
$ws.Cells.Item(189,7).Value = "후"
$ws.Cells.Item(189,8).Value = 260103
$ws.Cells.Item(189,9).Value = "느그클럽"
$ws.Cells.Item(189,10).Value = "공격"
$ws.Cells.Item(190,1).Value = "오공, 에이스, 아라곤"
$ws.Cells.Item(190,2).Value = "루"
$ws.Cells.Item(190,4).Value = "카일, 파이, 아멜리아"
$ws.Cells.Item(190,5).Value = "이린"
$ws.Cells.Item(190,6).Value = "오2아2"
$ws.Cells.Item(190,3).Value = "아1카1카2"
$ws.Cells.Item(190,7).Value = "후"
$ws.Cells.Item(190,8).Value = 260103
$ws.Cells.Item(190,9).Value = "느그클럽"
$ws.Cells.Item(190,10).Value = "공격"
$ws.Cells.Item(191,1).Value = "카일, 카구라, 파이"
$ws.Cells.Item(191,2).Value = "윈디"
$ws.Cells.Item(191,4).Value = "트루드, 카구라, 엘리시아"
$ws.Cells.Item(191,5).Value = "크리"
$ws.Cells.Item(191,6).Value = "엘1트2구2"
$ws.Cells.Item(191,3).Value = "구2카1카2"
$ws.Cells.Item(191,7).Value = "선"
$ws.Cells.Item(191,8).Value = 260103
$ws.Cells.Item(191,9).Value = "느그클럽"
$ws.Cells.Item(191,10).Value = "공격"
$ws.Cells.Item(192,1).Value = "플라튼, 콜트, 아멜리아"
$ws.Cells.Item(192,2).Value = "루"
$ws.Cells.Item(192,4).Value = "카일, 카구라, 아멜리아"
$ws.Cells.Item(192,5).Value = "이린"
$ws.Cells.Item(192,6).Value = "아2카1카2"
$ws.Cells.Item(192,3).Value = "플2콜1콜2"
$ws.Cells.Item(192,7).Value = "후"
$ws.Cells.Item(192,8).Value = 260103
$ws.Cells.Item(192,9).Value = "느그클럽"
$ws.Cells.Item(192,10).Value = "공격"
$ws.Cells.Item(193,1).Value = "연희, 로지, 키리엘"
$ws.Cells.Item(193,2).Value = "크리"
$ws.Cells.Item(193,4).Value = "플라튼, 스파이크, 엘리스"
$ws.Cells.Item(193,5).Value = "파이크"
$ws.Cells.Item(193,6).Value = "스2플2스1"
$ws.Cells.Item(193,3).Value = "키2연2연1"
$ws.Cells.Item(193,7).Value = "후"
$ws.Cells.Item(193,8).Value = 260103
$ws.Cells.Item(193,9).Value = "느그클럽"
$ws.Cells.Item(193,10).Value = "공격"
$ws.Cells.Item(194,1).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(194,2).Value = "맬패로"
$ws.Cells.Item(194,4).Value = "플라튼, 발리스타, 콜트"
$ws.Cells.Item(194,5).Value = "루"
$ws.Cells.Item(194,6).Value = "콜1콜2발2"
$ws.Cells.Item(194,3).Value = "바1프2"
$ws.Cells.Item(194,7).Value = "선"
$ws.Cells.Item(194,8).Value = 260103
$ws.Cells.Item(194,9).Value = "느그클럽"
$ws.Cells.Item(194,10).Value = "공격"
$ws.Cells.Item(195,1).Value = "플라튼, 콜트, 키리엘"
$ws.Cells.Item(195,2).Value = "루"
$ws.Cells.Item(195,4).Value = "카일, 카구라, 파이"
$ws.Cells.Item(195,5).Value = "이린"
$ws.Cells.Item(195,6).Value = "구2카1카2"
$ws.Cells.Item(195,3).Value = "키2콜1콜2"
$ws.Cells.Item(195,7).Value = "선"
$ws.Cells.Item(195,8).Value = 260103
$ws.Cells.Item(195,9).Value = "느그클럽"
$ws.Cells.Item(195,10).Value = "공격"
$ws.Cells.Item(196,1).Value = "연희, 바네사, 로지"
$ws.Cells.Item(196,2).Value = "노트"
$ws.Cells.Item(196,4).Value = "스파이크, 엘리스, 아라곤"
$ws.Cells.Item(196,5).Value = "파라곤"
$ws.Cells.Item(196,6).Value = "아2"
$ws.Cells.Item(196,3).Value = "바1연2연1"
$ws.Cells.Item(196,7).Value = "후"
$ws.Cells.Item(196,8).Value = 260103
$ws.Cells.Item(196,9).Value = "느그클럽"
$ws.Cells.Item(196,10).Value = "공격"
$ws.Cells.Item(197,1).Value = "트루드, 스파이크, 아멜리아"
$ws.Cells.Item(197,2).Value = "이린"
$ws.Cells.Item(197,4).Value = "밀리아, 바네사, 프레이야"
$ws.Cells.Item(197,5).Value = "유"
$ws.Cells.Item(197,6).Value = "밀2바1프2"
$ws.Cells.Item(197,3).Value = "아2트2아1"
$ws.Cells.Item(197,7).Value = "선"
$ws.Cells.Item(197,8).Value = 260103
$ws.Cells.Item(197,9).Value = "느그클럽"
$ws.Cells.Item(197,10).Value = "공격"
$ws.Cells.Item(198,1).Value = "콜트, 엘리시아, 델론즈"
$ws.Cells.Item(198,2).Value = "세리"
$ws.Cells.Item(198,4).Value = "카일, 카구라, 파이"
$ws.Cells.Item(198,5).Value = "이린"
$ws.Cells.Item(198,6).Value = "콜2엘2엘1"
$ws.Cells.Item(198,3).Value = "콜2엘2엘1"
$ws.Cells.Item(198,7).Value = "후"
$ws.Cells.Item(198,8).Value = 260103
$ws.Cells.Item(198,9).Value = "느그클럽"
$ws.Cells.Item(198,10).Value = "공격"
$ws.Cells.Item(199,1).Value = "밀리아, 쥬리, 멜키르"
$ws.Cells.Item(199,2).Value = "루"
$ws.Cells.Item(199,4).Value = "밀리아, 오공, 프레이야"
$ws.Cells.Item(199,5).Value = "연지"
$ws.Cells.Item(199,6).Value = "프2프1오2"
$ws.Cells.Item(199,3).Value = "밀2멜2쥬2"
$ws.Cells.Item(199,7).Value = "선"
$ws.Cells.Item(199,8).Value = 260103
$ws.Cells.Item(199,9).Value = "느그클럽"
$ws.Cells.Item(199,10).Value = "공격"
$ws.Cells.Item(200,1).Value = "플라튼, 실베스타, 아멜리아"
$ws.Cells.Item(200,2).Value = "파이크"
$ws.Cells.Item(200,4).Value = "플라튼, 트루드, 아멜리아"
$ws.Cells.Item(200,5).Value = "파이크"
$ws.Cells.Item(200,6).Value = "아2아1트2"
$ws.Cells.Item(200,3).Value = "아2실2실1"
$ws.Cells.Item(200,7).Value = "후"
$ws.Cells.Item(200,8).Value = 260103
$ws.Cells.Item(200,9).Value = "느그클럽"
$ws.Cells.Item(200,10).Value = "공격"
$ws.Cells.Item(201,1).Value = "프레이야, 바네사, 아멜리아"
$ws.Cells.Item(201,2).Value = "노트"
$ws.Cells.Item(201,4).Value = "프레이야, 바네사, 밀리아"
$ws.Cells.Item(201,5).Value = "맬패로"
$ws.Cells.Item(201,6).Value = "바1프2프1"
$ws.Cells.Item(201,3).Value = "바1프2바2"
$ws.Cells.Item(201,7).Value = "선"
$ws.Cells.Item(201,8).Value = 260103
$ws.Cells.Item(201,9).Value = "느그클럽"
$ws.Cells.Item(201,10).Value = "공격"
$ws.Cells.Item(202,1).Value = "트루드, 크리스, 엘리시아"
$ws.Cells.Item(202,2).Value = "카람"
$ws.Cells.Item(202,4).Value = "겔리두스, 스파이크, 크리스"
$ws.Cells.Item(202,5).Value = "루"
$ws.Cells.Item(202,6).Value = "겔2스2크2"
$ws.Cells.Item(202,3).Value = "크1엘2트2"
$ws.Cells.Item(202,7).Value = "후"
$ws.Cells.Item(202,8).Value = 260103
$ws.Cells.Item(202,9).Value = "느그클럽"
$ws.Cells.Item(202,10).Value = "공격"
$ws.Cells.Item(203,1).Value = "프레이야, 키리엘, 쥬리"
$ws.Cells.Item(203,2).Value = "크리"
$ws.Cells.Item(203,4).Value = "카일, 카구라, 콜트"
$ws.Cells.Item(203,5).Value = "이린"
$ws.Cells.Item(203,6).Value = "구2카1카2"
$ws.Cells.Item(203,3).Value = "쥬2프2프1"
$ws.Cells.Item(203,7).Value = "후"
$ws.Cells.Item(203,8).Value = 260103
$ws.Cells.Item(203,9).Value = "느그클럽"
$ws.Cells.Item(203,10).Value = "공격"
$ws.Cells.Item(204,1).Value = "트루드, 엘리시아, 콜트"
$ws.Cells.Item(204,2).Value = "이린"
$ws.Cells.Item(204,4).Value = "스파이크, 엘리스, 트루드"
$ws.Cells.Item(204,5).Value = "루"
$ws.Cells.Item(204,6).Value = "엘1트2스2"
$ws.Cells.Item(204,3).Value = "콜1콜2"
$ws.Cells.Item(204,7).Value = "후"
$ws.Cells.Item(204,8).Value = 260103
$ws.Cells.Item(204,9).Value = "느그클럽"
$ws.Cells.Item(204,10).Value = "공격"
$ws.Cells.Item(205,1).Value = "밀리아, 연희, 린"
$ws.Cells.Item(205,2).Value = "유"
$ws.Cells.Item(205,4).Value = "오공, 엘리시아, 겔리두스"
$ws.Cells.Item(205,5).Value = "파라곤"
$ws.Cells.Item(205,6).Value = "오2겔1겔2"
$ws.Cells.Item(205,3).Value = "밀2연2밀1"
$ws.Cells.Item(205,7).Value = "후"
$ws.Cells.Item(205,8).Value = 260103
$ws.Cells.Item(205,9).Value = "느그클럽"
$ws.Cells.Item(205,10).Value = "공격"
$ws.Cells.Item(206,1).Value = "연희, 바네사, 로지"
$ws.Cells.Item(206,2).Value = "연지"
$ws.Cells.Item(206,4).Value = "프레이야, 트루드, 콜트"
$ws.Cells.Item(206,5).Value = "카람"
$ws.Cells.Item(206,6).Value = "콜1콜2프1"
$ws.Cells.Item(206,3).Value = "바1연2바2"
$ws.Cells.Item(206,7).Value = "후"
$ws.Cells.Item(206,8).Value = 260103
$ws.Cells.Item(206,9).Value = "느그클럽"
$ws.Cells.Item(206,10).Value = "공격"
$ws.Cells.Item(207,1).Value = "연희, 키리엘, 콜트"
$ws.Cells.Item(207,2).Value = "델로"
$ws.Cells.Item(207,4).Value = "프레이야, 바네사, 밀리아"
$ws.Cells.Item(207,5).Value = "델로"
$ws.Cells.Item(207,6).Value = "밀2프2프1"
$ws.Cells.Item(207,3).Value = "연2콜2키2"
$ws.Cells.Item(207,7).Value = "후"
$ws.Cells.Item(207,8).Value = 260103
$ws.Cells.Item(207,9).Value = "느그클럽"
$ws.Cells.Item(207,10).Value = "공격"

# Update the sheet view: active selection to match the state the workbook
# was left in after the edit (scrolled down to the newly added rows).
$ws.Activate()
$ws.Range("A160").Select()
$ws.Range("O199").Select()

Write-Output "Added 30 new rows (178-207) to Sheet1."
